$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 37.44676833333333
$ws.Cells.Item(2, 8).Value = 112.340305
$ws.Cells.Item(2, 9).Value = 0.6029245723174422
$ws.Cells.Item(2, 10).Value = 0.6029245723174423
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.1474273333333333
$ws.Cells.Item(2, 14).Value = 0.442282
$ws.Cells.Item(2, 15).Value = 0.1588601259223368
$ws.Cells.Item(2, 16).Value = 0.1588601259223368
$ws.Cells.Item(2, 17).Value = 5.520677197334444
$ws.Cells.Item(2, 18).Value = 49.68609477601
$ws.Cells.Item(2, 19).Value = 0.09578067348001991
$ws.Cells.Item(2, 20).Value = 0.09578067348001992

$ws.Cells.Item(3, 7).Value = 37.44676833333333
$ws.Cells.Item(3, 8).Value = 112.340305
$ws.Cells.Item(3, 9).Value = 0.6029245723174422
$ws.Cells.Item(3, 10).Value = 0.6029245723174423
$ws.Cells.Item(3, 15).Value = 0.4626735347223893
$ws.Cells.Item(3, 16).Value = 0.4626735347223893
$ws.Cells.Item(3, 17).Value = 16.07874359989333
$ws.Cells.Item(3, 18).Value = 144.70869239904
$ws.Cells.Item(3, 19).Value = 0.2789572430450958
$ws.Cells.Item(3, 20).Value = 0.2789572430450959

$ws.Cells.Item(4, 7).Value = 37.44676833333333
$ws.Cells.Item(4, 8).Value = 112.340305
$ws.Cells.Item(4, 9).Value = 0.6029245723174422
$ws.Cells.Item(4, 10).Value = 0.6029245723174423
$ws.Cells.Item(4, 13).Value = 0.01780266666666666
$ws.Cells.Item(4, 14).Value = 0.053408
$ws.Cells.Item(4, 15).Value = 0.01918323966442261
$ws.Cells.Item(4, 16).Value = 0.01918323966442261
$ws.Cells.Item(4, 17).Value = 0.6666523343822222
$ws.Cells.Item(4, 18).Value = 5.99987100944
$ws.Cells.Item(4, 19).Value = 0.01156604657033499
$ws.Cells.Item(4, 20).Value = 0.011566046570335

$ws.Cells.Item(5, 7).Value = 37.44676833333333
$ws.Cells.Item(5, 8).Value = 112.340305
$ws.Cells.Item(5, 9).Value = 0.6029245723174422
$ws.Cells.Item(5, 10).Value = 0.6029245723174423
$ws.Cells.Item(5, 13).Value = 0.3334263333333333
$ws.Cells.Item(5, 14).Value = 1.000279
$ws.Cells.Item(5, 15).Value = 0.3592830996908513
$ws.Cells.Item(5, 16).Value = 0.3592830996908513
$ws.Cells.Item(5, 17).Value = 12.48573866056611
$ws.Cells.Item(5, 18).Value = 112.371647945095
$ws.Cells.Item(5, 19).Value = 0.2166206092219914
$ws.Cells.Item(5, 20).Value = 0.2166206092219915

$ws.Cells.Item(6, 9).Value = 0.1838793176915316
$ws.Cells.Item(6, 10).Value = 0.1838793176915316
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.1474273333333333
$ws.Cells.Item(6, 14).Value = 0.442282
$ws.Cells.Item(6, 15).Value = 0.1588601259223368
$ws.Cells.Item(6, 16).Value = 0.1588601259223368
$ws.Cells.Item(6, 17).Value = 1.683690469504667
$ws.Cells.Item(6, 18).Value = 15.153214225542
$ws.Cells.Item(6, 19).Value = 0.02921109156299008
$ws.Cells.Item(6, 20).Value = 0.02921109156299008

$ws.Cells.Item(7, 9).Value = 0.1838793176915316
$ws.Cells.Item(7, 10).Value = 0.1838793176915316
$ws.Cells.Item(7, 15).Value = 0.4626735347223893
$ws.Cells.Item(7, 16).Value = 0.4626735347223893
$ws.Cells.Item(7, 19).Value = 0.08507609387868212
$ws.Cells.Item(7, 20).Value = 0.08507609387868212

$ws.Cells.Item(8, 9).Value = 0.1838793176915316
$ws.Cells.Item(8, 10).Value = 0.1838793176915316
$ws.Cells.Item(8, 13).Value = 0.01780266666666666
$ws.Cells.Item(8, 14).Value = 0.053408
$ws.Cells.Item(8, 15).Value = 0.01918323966442261
$ws.Cells.Item(8, 16).Value = 0.01918323966442261
$ws.Cells.Item(8, 17).Value = 0.2033149452053333
$ws.Cells.Item(8, 18).Value = 1.829834506848
$ws.Cells.Item(8, 19).Value = 0.003527401020607155
$ws.Cells.Item(8, 20).Value = 0.003527401020607156

$ws.Cells.Item(9, 9).Value = 0.1838793176915316
$ws.Cells.Item(9, 10).Value = 0.1838793176915316
$ws.Cells.Item(9, 13).Value = 0.3334263333333333
$ws.Cells.Item(9, 14).Value = 1.000279
$ws.Cells.Item(9, 15).Value = 0.3592830996908513
$ws.Cells.Item(9, 16).Value = 0.3592830996908513
$ws.Cells.Item(9, 17).Value = 3.807887771027667
$ws.Cells.Item(9, 18).Value = 34.270989939249
$ws.Cells.Item(9, 19).Value = 0.06606473122925227
$ws.Cells.Item(9, 20).Value = 0.06606473122925227

$ws.Cells.Item(10, 7).Value = 1.244612333333333
$ws.Cells.Item(10, 8).Value = 3.733837
$ws.Cells.Item(10, 9).Value = 0.02003930892236799
$ws.Cells.Item(10, 10).Value = 0.02003930892236799
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.1474273333333333
$ws.Cells.Item(10, 14).Value = 0.442282
$ws.Cells.Item(10, 15).Value = 0.1588601259223368
$ws.Cells.Item(10, 16).Value = 0.1588601259223368
$ws.Cells.Item(10, 17).Value = 0.1834898773371111
$ws.Cells.Item(10, 18).Value = 1.651408896034
$ws.Cells.Item(10, 19).Value = 0.003183447138803986
$ws.Cells.Item(10, 20).Value = 0.003183447138803986

$ws.Cells.Item(11, 7).Value = 1.244612333333333
$ws.Cells.Item(11, 8).Value = 3.733837
$ws.Cells.Item(11, 9).Value = 0.02003930892236799
$ws.Cells.Item(11, 10).Value = 0.02003930892236799
$ws.Cells.Item(11, 15).Value = 0.4626735347223893
$ws.Cells.Item(11, 16).Value = 0.4626735347223893
$ws.Cells.Item(11, 17).Value = 0.5344066652373333
$ws.Cells.Item(11, 18).Value = 4.809659987136
$ws.Cells.Item(11, 19).Value = 0.009271657892505914
$ws.Cells.Item(11, 20).Value = 0.009271657892505914

$ws.Cells.Item(12, 7).Value = 1.244612333333333
$ws.Cells.Item(12, 8).Value = 3.733837
$ws.Cells.Item(12, 9).Value = 0.02003930892236799
$ws.Cells.Item(12, 10).Value = 0.02003930892236799
$ws.Cells.Item(12, 13).Value = 0.01780266666666666
$ws.Cells.Item(12, 14).Value = 0.053408
$ws.Cells.Item(12, 15).Value = 0.01918323966442261
$ws.Cells.Item(12, 16).Value = 0.01918323966442261
$ws.Cells.Item(12, 17).Value = 0.02215741849955555
$ws.Cells.Item(12, 18).Value = 0.199416766496
$ws.Cells.Item(12, 19).Value = 0.0003844188657671876
$ws.Cells.Item(12, 20).Value = 0.0003844188657671876

$ws.Cells.Item(13, 7).Value = 1.244612333333333
$ws.Cells.Item(13, 8).Value = 3.733837
$ws.Cells.Item(13, 9).Value = 0.02003930892236799
$ws.Cells.Item(13, 10).Value = 0.02003930892236799
$ws.Cells.Item(13, 13).Value = 0.3334263333333333
$ws.Cells.Item(13, 14).Value = 1.000279
$ws.Cells.Item(13, 15).Value = 0.3592830996908513
$ws.Cells.Item(13, 16).Value = 0.3592830996908513
$ws.Cells.Item(13, 17).Value = 0.4149865267247777
$ws.Cells.Item(13, 18).Value = 3.734878740523
$ws.Cells.Item(13, 19).Value = 0.007199785025290905
$ws.Cells.Item(13, 20).Value = 0.007199785025290905

$ws.Cells.Item(14, 7).Value = 10.03858766666667
$ws.Cells.Item(14, 8).Value = 30.115763
$ws.Cells.Item(14, 9).Value = 0.1616297332180864
$ws.Cells.Item(14, 10).Value = 0.1616297332180864
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.1474273333333333
$ws.Cells.Item(14, 14).Value = 0.442282
$ws.Cells.Item(14, 15).Value = 0.1588601259223368
$ws.Cells.Item(14, 16).Value = 0.1588601259223368
$ws.Cells.Item(14, 17).Value = 1.479962210129556
$ws.Cells.Item(14, 18).Value = 13.319659891166
$ws.Cells.Item(14, 19).Value = 0.02567651977181889
$ws.Cells.Item(14, 20).Value = 0.02567651977181889

$ws.Cells.Item(15, 7).Value = 10.03858766666667
$ws.Cells.Item(15, 8).Value = 30.115763
$ws.Cells.Item(15, 9).Value = 0.1616297332180864
$ws.Cells.Item(15, 10).Value = 0.1616297332180864
$ws.Cells.Item(15, 15).Value = 0.4626735347223893
$ws.Cells.Item(15, 16).Value = 0.4626735347223893
$ws.Cells.Item(15, 17).Value = 4.310328617962666
$ws.Cells.Item(15, 18).Value = 38.792957561664
$ws.Cells.Item(15, 19).Value = 0.0747817999842488
$ws.Cells.Item(15, 20).Value = 0.0747817999842488

$ws.Cells.Item(16, 7).Value = 10.03858766666667
$ws.Cells.Item(16, 8).Value = 30.115763
$ws.Cells.Item(16, 9).Value = 0.1616297332180864
$ws.Cells.Item(16, 10).Value = 0.1616297332180864
$ws.Cells.Item(16, 13).Value = 0.01780266666666666
$ws.Cells.Item(16, 14).Value = 0.053408
$ws.Cells.Item(16, 15).Value = 0.01918323966442261
$ws.Cells.Item(16, 16).Value = 0.01918323966442261
$ws.Cells.Item(16, 17).Value = 0.1787136300337777
$ws.Cells.Item(16, 18).Value = 1.608422670304
$ws.Cells.Item(16, 19).Value = 0.003100581909219238
$ws.Cells.Item(16, 20).Value = 0.003100581909219239

$ws.Cells.Item(17, 7).Value = 10.03858766666667
$ws.Cells.Item(17, 8).Value = 30.115763
$ws.Cells.Item(17, 9).Value = 0.1616297332180864
$ws.Cells.Item(17, 10).Value = 0.1616297332180864
$ws.Cells.Item(17, 13).Value = 0.3334263333333333
$ws.Cells.Item(17, 14).Value = 1.000279
$ws.Cells.Item(17, 15).Value = 0.3592830996908513
$ws.Cells.Item(17, 16).Value = 0.3592830996908513
$ws.Cells.Item(17, 17).Value = 3.347129477541889
$ws.Cells.Item(17, 18).Value = 30.124165297877
$ws.Cells.Item(17, 19).Value = 0.05807083155279941
$ws.Cells.Item(17, 20).Value = 0.05807083155279941

$ws.Cells.Item(18, 7).Value = 1.327177333333333
$ws.Cells.Item(18, 8).Value = 3.981532
$ws.Cells.Item(18, 9).Value = 0.02136867510078605
$ws.Cells.Item(18, 10).Value = 0.02136867510078605
$ws.Cells.Item(18, 11).Value = 2
$ws.Cells.Item(18, 12).Value = 0.6666666666666666
$ws.Cells.Item(18, 13).Value = 0.1474273333333333
$ws.Cells.Item(18, 14).Value = 0.442282
$ws.Cells.Item(18, 15).Value = 0.1588601259223368
$ws.Cells.Item(18, 16).Value = 0.1588601259223368
$ws.Cells.Item(18, 17).Value = 0.1956622151137778
$ws.Cells.Item(18, 18).Value = 1.760959936024
$ws.Cells.Item(18, 19).Value = 0.003394630417304374
$ws.Cells.Item(18, 20).Value = 0.003394630417304374

$ws.Cells.Item(19, 7).Value = 1.327177333333333
$ws.Cells.Item(19, 8).Value = 3.981532
$ws.Cells.Item(19, 9).Value = 0.02136867510078605
$ws.Cells.Item(19, 10).Value = 0.02136867510078605
$ws.Cells.Item(19, 15).Value = 0.4626735347223893
$ws.Cells.Item(19, 16).Value = 0.4626735347223893
$ws.Cells.Item(19, 17).Value = 0.5698580946773334
$ws.Cells.Item(19, 18).Value = 5.128722852096
$ws.Cells.Item(19, 19).Value = 0.009886720441214992
$ws.Cells.Item(19, 20).Value = 0.00988672044121499

$ws.Cells.Item(20, 7).Value = 1.327177333333333
$ws.Cells.Item(20, 8).Value = 3.981532
$ws.Cells.Item(20, 9).Value = 0.02136867510078605
$ws.Cells.Item(20, 10).Value = 0.02136867510078605
$ws.Cells.Item(20, 13).Value = 0.01780266666666666
$ws.Cells.Item(20, 14).Value = 0.053408
$ws.Cells.Item(20, 15).Value = 0.01918323966442261
$ws.Cells.Item(20, 16).Value = 0.01918323966442261
$ws.Cells.Item(20, 17).Value = 0.02362729567288889
$ws.Cells.Item(20, 18).Value = 0.212645661056
$ws.Cells.Item(20, 19).Value = 0.0004099204157695588
$ws.Cells.Item(20, 20).Value = 0.0004099204157695588

$ws.Cells.Item(21, 7).Value = 1.327177333333333
$ws.Cells.Item(21, 8).Value = 3.981532
$ws.Cells.Item(21, 9).Value = 0.02136867510078605
$ws.Cells.Item(21, 10).Value = 0.02136867510078605
$ws.Cells.Item(21, 13).Value = 0.3334263333333333
$ws.Cells.Item(21, 14).Value = 1.000279
$ws.Cells.Item(21, 15).Value = 0.3592830996908513
$ws.Cells.Item(21, 16).Value = 0.3592830996908513
$ws.Cells.Item(21, 17).Value = 0.4425158719364445
$ws.Cells.Item(21, 18).Value = 3.982642847428
$ws.Cells.Item(21, 19).Value = 0.007677403826497127
$ws.Cells.Item(21, 20).Value = 0.007677403826497125

$ws.Cells.Item(22, 7).Value = 0.630923
$ws.Cells.Item(22, 8).Value = 1.892769
$ws.Cells.Item(22, 9).Value = 0.01015839274978569
$ws.Cells.Item(22, 10).Value = 0.01015839274978569
$ws.Cells.Item(22, 11).Value = 2
$ws.Cells.Item(22, 12).Value = 0.6666666666666666
$ws.Cells.Item(22, 13).Value = 0.1474273333333333
$ws.Cells.Item(22, 14).Value = 0.442282
$ws.Cells.Item(22, 15).Value = 0.1588601259223368
$ws.Cells.Item(22, 16).Value = 0.1588601259223368
$ws.Cells.Item(22, 17).Value = 0.09301529542866667
$ws.Cells.Item(22, 18).Value = 0.837137658858
$ws.Cells.Item(22, 19).Value = 0.001613763551399507
$ws.Cells.Item(22, 20).Value = 0.001613763551399507

$ws.Cells.Item(23, 7).Value = 0.630923
$ws.Cells.Item(23, 8).Value = 1.892769
$ws.Cells.Item(23, 9).Value = 0.01015839274978569
$ws.Cells.Item(23, 10).Value = 0.01015839274978569
$ws.Cells.Item(23, 15).Value = 0.4626735347223893
$ws.Cells.Item(23, 16).Value = 0.4626735347223893
$ws.Cells.Item(23, 17).Value = 0.270903194048
$ws.Cells.Item(23, 18).Value = 2.438128746432
$ws.Cells.Item(23, 19).Value = 0.004700019480641636
$ws.Cells.Item(23, 20).Value = 0.004700019480641636

$ws.Cells.Item(24, 7).Value = 0.630923
$ws.Cells.Item(24, 8).Value = 1.892769
$ws.Cells.Item(24, 9).Value = 0.01015839274978569
$ws.Cells.Item(24, 10).Value = 0.01015839274978569
$ws.Cells.Item(24, 13).Value = 0.01780266666666666
$ws.Cells.Item(24, 14).Value = 0.053408
$ws.Cells.Item(24, 15).Value = 0.01918323966442261
$ws.Cells.Item(24, 16).Value = 0.01918323966442261
$ws.Cells.Item(24, 17).Value = 0.01123211186133333
$ws.Cells.Item(24, 18).Value = 0.101089006752
$ws.Cells.Item(24, 19).Value = 0.0001948708827244719
$ws.Cells.Item(24, 20).Value = 0.0001948708827244719

$ws.Cells.Item(25, 7).Value = 0.630923
$ws.Cells.Item(25, 8).Value = 1.892769
$ws.Cells.Item(25, 9).Value = 0.01015839274978569
$ws.Cells.Item(25, 10).Value = 0.01015839274978569
$ws.Cells.Item(25, 13).Value = 0.3334263333333333
$ws.Cells.Item(25, 14).Value = 1.000279
$ws.Cells.Item(25, 15).Value = 0.3592830996908513
$ws.Cells.Item(25, 16).Value = 0.3592830996908513
$ws.Cells.Item(25, 17).Value = 0.2103663425056667
$ws.Cells.Item(25, 18).Value = 1.893297082551
$ws.Cells.Item(25, 19).Value = 0.003649738835020072
$ws.Cells.Item(25, 20).Value = 0.003649738835020072
